$d = $word.ActiveDocument

$replacements = @(
    @("457×4=", "553×5="),
    @("988×5=", "568×4="),
    @("393×7=", "404×6="),
    @("157×4=", "650×8="),
    @("357×2=", "151×4="),
    @("143×8=", "811×7="),
    @("916×5=", "167×7="),
    @("461×4=", "488×5="),
    @("458×5=", "810×7="),
    @("616×5=", "603×2="),
    @("859×6=", "900×4="),
    @("329×2=", "178×2="),
    @("425×2=", "997×3="),
    @("216×8=", "222×3="),
    @("245×3=", "937×2="),
    @("276×7=", "697×8="),
    @("344×6=", "367×5="),
    @("136×7=", "788×8="),
    @("516×8=", "401×4="),
    @("337×3=", "187×8="),
    @("188×5=", "263×4="),
    @("136×4=", "285×7="),
    @("972×8=", "504×2="),
    @("920×7=", "232×7="),
    @("526×7=", "806×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
